$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.756.24'
$ws.Range('E2').Value = '  +4.93%  '
$ws.Range('D3').Value = '2.759.64'
$ws.Range('E3').Value = '  +4.03%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.49'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.57'
$ws.Range('E6').Value = '  +6.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  +1.94%  '
$ws.Range('D9').Value = '2.757.47'
$ws.Range('E9').Value = '  +3.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.70'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('E11').Value = '  +5.41%  '
$ws.Range('E12').Value = '  +4.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.390'
$ws.Range('E13').Value = '  +3.96%  '
$ws.Range('D14').Value = '3.239.77'
$ws.Range('E14').Value = '  +3.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.42'
$ws.Range('E15').Value = '  +2.95%  '
$ws.Range('D16').Value = '63.671.54'
$ws.Range('E16').Value = '  +4.86%  '
$ws.Range('E17').Value = '  +6.54%  '
$ws.Range('D18').Value = '2.749.90'
$ws.Range('E18').Value = '  +3.39%  '
$ws.Range('E19').Value = '  +3.69%  '
$ws.Range('E20').Value = '  +3.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '360.62'
$ws.Range('E21').Value = '  +2.91%  '
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '66.09'
$ws.Range('E25').Value = '  +3.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.170'
$ws.Range('E26').Value = '  +5.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.55'
$ws.Range('E27').Value = '  +4.90%  '
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('D29').Value = '0.0₃0913'
$ws.Range('E29').Value = '  +12.74%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('E31').Value = '  +5.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '171.49'
$ws.Range('E32').Value = '  +2.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.22'
$ws.Range('E33').Value = '  +15.27%  '
$ws.Range('E34').Value = '  -0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.50'
$ws.Range('E35').Value = '  +3.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.78'
$ws.Range('E36').Value = '  +7.67%  '
$ws.Range('E37').Value = '  +9.12%  '
$ws.Range('E38').Value = '  +10.38%  '
$ws.Range('E39').Value = '  +14.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '344.15'
$ws.Range('E40').Value = '  +5.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.23'
$ws.Range('E41').Value = '  +5.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.40'
$ws.Range('E42').Value = '  +2.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.61'
$ws.Range('E43').Value = '  +6.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.72'
$ws.Range('E44').Value = '  +5.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.85'
$ws.Range('E45').Value = '  +6.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '140.08'
$ws.Range('E46').Value = '  +3.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0592'
$ws.Range('E47').Value = '  +6.21%  '
$ws.Range('E48').Value = '  +5.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0255'
$ws.Range('E49').Value = '  +3.65%  '
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.997'
$ws.Range('E51').Value = '  -0.15%  '
